$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before row 654, pushing the existing rows 654-668 down to 657-671.
$ws.Rows("654:656").Insert()

# Common values shared by every Brócoli record in this block.
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$categoriaId = 100112023
$categoria = "Brócoli"
$variedad = "Sin especificar"
$unidad = "$/unidad"
$origen = "Región de Arica y Parinacota"
$kgUnidades = 1
$clasificacion = "Hortaliza"

function Set-BrocoliRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $calidad
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $precioMin
    $ws.Cells.Item($row, 12).Value = $precioMax
    $ws.Cells.Item($row, 13).Value = $precioProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $precioKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}

Set-BrocoliRow 654 45239 "Primera" 200 600 700 650 650
Set-BrocoliRow 655 45239 "Segunda" 800 500 600 544 544
Set-BrocoliRow 656 45239 "Tercera" 480 300 400 342 342
